$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 188.83333
$ws.Range("I38").Value = 188.83333
$ws.Range("K38").Value = 566.49999
$ws.Range("M38").Value = -194.49999

$ws.Range("H62").Value = 7624.875
$ws.Range("I62").Value = 4999
$ws.Range("K62").Value = 4999
$ws.Range("M62").Value = -4375

$ws.Range("H65").Value = 7624.875
$ws.Range("I65").Value = 4999
$ws.Range("K65").Value = 24995
$ws.Range("M65").Value = -21875

$ws.Range("H86").Value = 4619.7856
$ws.Range("I86").Value = 3399.25
$ws.Range("K86").Value = 3399.25
$ws.Range("M86").Value = -2276.25

$ws.Range("H89").Value = 4619.7856
$ws.Range("I89").Value = 3399.25
$ws.Range("K89").Value = 16996.25
$ws.Range("M89").Value = -11380.25

$ws.Range("H113").Value = 2999.5
$ws.Range("I113").Value = 2999.5
$ws.Range("K113").Value = 2999.5
$ws.Range("M113").Value = 254.5

$ws.Range("H116").Value = 5732.421
$ws.Range("J116").Value = 6504.5557
$ws.Range("L116").Value = 6504.5557
$ws.Range("N116").Value = -13388.5557

$ws.Range("H132").Value = 1286.3214
$ws.Range("I132").Value = 1319.3704
$ws.Range("K132").Value = 3958.1112
$ws.Range("M132").Value = -1428.1112

$ws.Range("H137").Value = 2921.647
$ws.Range("I137").Value = 2089.8462
$ws.Range("J137").Value = 5625
$ws.Range("K137").Value = 6269.5386
$ws.Range("L137").Value = 16875
$ws.Range("M137").Value = -3719.5386
$ws.Range("N137").Value = -21975

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1236.5
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H17").Value = 44
$ws.Range("I17").Value = 7.3333335
$ws.Range("K17").Value = 7.3333335
$ws.Range("M17").Value = 165.6666665

$ws.Range("H32").Value = 7439.643
$ws.Range("I32").Value = 5591.184
$ws.Range("K32").Value = 5591.184
$ws.Range("M32").Value = -5304.184

$ws.Range("H63").Value = 6537.5625
$ws.Range("I63").Value = 5010.5557
$ws.Range("J63").Value = 8500.857
$ws.Range("K63").Value = 5010.5557
$ws.Range("L63").Value = 8500.857
$ws.Range("M63").Value = -4324.5557
$ws.Range("N63").Value = -9872.857

$ws.Range("H66").Value = 6537.5625
$ws.Range("I66").Value = 5010.5557
$ws.Range("J66").Value = 8500.857
$ws.Range("K66").Value = 25052.7785
$ws.Range("L66").Value = 42504.285
$ws.Range("M66").Value = -21620.7785
$ws.Range("N66").Value = -49368.285

$ws.Range("H116").Value = 1236.5
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1236.5
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H22").Value = 626
$ws.Range("I22").Value = 626
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 626
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -453

$ws.Range("H88").Value = 19499.834
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 21399.8
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 21399.8
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -22211.8

$ws.Range("H91").Value = 19499.834
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 21399.8
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 21399.8
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -24207.8

$ws.Range("H99").Value = 4680.5
$ws.Range("I99").Value = 4593
$ws.Range("J99").Value = 4855.5
$ws.Range("K99").Value = 4593
$ws.Range("L99").Value = 4855.5
$ws.Range("M99").Value = -3095
$ws.Range("N99").Value = -7851.5

$ws.Range("H105").Value = 4095.6365
$ws.Range("J105").Value = 15499.5
$ws.Range("L105").Value = 15499.5
$ws.Range("N105").Value = -18993.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 392.2
$ws.Range("I23").Value = 467
$ws.Range("J23").Value = 280
$ws.Range("K23").Value = 1401
$ws.Range("L23").Value = 840
$ws.Range("M23").Value = -1166
$ws.Range("N23").Value = -1310

$ws.Range("H38").Value = 162.14285
$ws.Range("I38").Value = 121.333336
$ws.Range("K38").Value = 364.000008
$ws.Range("M38").Value = -17.00000799999998

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H121").Value = 1200.909
$ws.Range("I121").Value = 715.4
$ws.Range("J121").Value = 1605.5
$ws.Range("K121").Value = 2146.2
$ws.Range("L121").Value = 4816.5
$ws.Range("M121").Value = -836.1999999999998
$ws.Range("N121").Value = -7436.5

$ws.Range("H131").Value = 1898
$ws.Range("J131").Value = 1915.1666
$ws.Range("L131").Value = 5745.4998
$ws.Range("N131").Value = -15825.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 4583.3335
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 4583.3335
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 4583.3335
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -8327.333500000001

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3344.923
$ws.Range("I100").Value = 1974.75
$ws.Range("K100").Value = 1974.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
